# Applies the "Updated cryptos list" data refresh to Sheet1 (columns D = Price, E = Volume(1h)).
# Numeric-looking Price values are written with a leading apostrophe so Excel keeps them as
# literal text (matching the workbook's existing inline-string convention) instead of converting
# them to real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.653.50'
$ws.Range("D3").Value = '2.239.93'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '''305.66'
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").Value = '''93.05'
$ws.Range("E6").Value = '  -1.84%  '
$ws.Range("D7").Value = '''0.568'
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("D10").Value = '''34.56'
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").Value = '''0.0797'
$ws.Range("E11").Value = '  -1.76%  '
$ws.Range("D12").Value = '''7.11'
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("D14").Value = '2.582.22'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '2.330.49'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("D17").Value = '''13.44'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '44.391.09'
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("D19").Value = '0.0₃0929'
$ws.Range("E19").Value = '  -3.65%  '
$ws.Range("E20").Value = '  -4.09%  '
$ws.Range("D21").Value = '''11.65'
$ws.Range("E21").Value = '  -4.39%  '
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").Value = '''236.93'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  +3.45%  '
$ws.Range("E28").Value = '  -2.18%  '
$ws.Range("D29").Value = '''36.81'
$ws.Range("E29").Value = '  -5.43%  '
$ws.Range("D30").Value = '''19.88'
$ws.Range("E30").Value = '  -0.81%  '
$ws.Range("D31").Value = '''5.79'
$ws.Range("E31").Value = '  -1.18%  '
$ws.Range("D32").Value = '''148.55'
$ws.Range("E32").Value = '  -3.18%  '
$ws.Range("E33").Value = '  +0.43%  '
$ws.Range("E34").Value = '  -2.69%  '
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("D38").Value = '''1.85'
$ws.Range("E38").Value = '  +4.74%  '
$ws.Range("D39").Value = '''14.96'
$ws.Range("E39").Value = '  +4.71%  '
$ws.Range("E40").Value = '  -4.85%  '
$ws.Range("E41").Value = '  -1.40%  '
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '1.810.18'
$ws.Range("E44").Value = '  +3.58%  '
$ws.Range("D45").Value = '''1.75'
$ws.Range("E45").Value = '  +9.21%  '
$ws.Range("D46").Value = '''81.24'
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("D47").Value = '''0.186'
$ws.Range("E47").Value = '  -3.32%  '
$ws.Range("D48").Value = '''97.27'
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("E49").Value = '  -2.38%  '
$ws.Range("D50").Value = '''68.54'
$ws.Range("E50").Value = '  +2.05%  '
$ws.Range("D51").Value = '''53.50'
$ws.Range("E51").Value = '  -2.28%  '

Write-Host "Updated cryptos list"
